$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$headers = @("NB", "SCH", "LENGTH", "PROCESS DONE", "NEXT PROCESS")
$cols = @(16, 17, 18, 19, 20)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $cols[$i])
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Font.ThemeColor = 0
    $cell.Interior.Color = 13140480
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = 2
}

$ws.Columns.Item(19).ColumnWidth = 14.5703125
$ws.Columns.Item(20).ColumnWidth = 14

$ws.Range("S4").Select()
